$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Amazon + 1P Sales Trend update — append weekly summary rows (units_ordered /
# ordered_product_sales) for a handful of SKUs to the bottom of the sheet.

$rows = @(
    @{ Row = 39; SKU = "FBA79612"; Model = "WM1ML";      Asin = "B0DP2WC5VW"; Units = 35; Sales = 82991.62 },
    @{ Row = 40; SKU = "FBA79613"; Model = "MS1ML";      Asin = "B0DP2VVRND"; Units = 47; Sales = 75200 },
    @{ Row = 41; SKU = "FBA79476"; Model = "WM-GS1M-BK"; Asin = "B0DB5VG39T"; Units = 26; Sales = 56142.32 },
    @{ Row = 42; SKU = "FBA79617"; Model = "HDWF1ML";    Asin = "B0DP32F346"; Units = 20; Sales = 50813.6 },
    @{ Row = 43; SKU = "FBA79616"; Model = "HD1ML";      Asin = "B0DP3194QN"; Units = 2;  Sales = 6438.98 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.SKU
    $ws.Range("B$n").Value = $r.Model
    $ws.Range("C$n").Value = $r.Asin
    $ws.Range("D$n").Value = $r.Asin
    $ws.Range("P$n").Value = $r.Units
    $ws.Range("T$n").Value = $r.Sales
}

# Leave the view scrolled to / focused on the newly added total so the sheet
# opens where the new data is, matching the saved workbook view state.
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("T39").Select()
